$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small wording tweak to the description of the last use case. Do this
# before inserting the new column / adding new text below, so the edit
# lands in-place on the existing shared string instead of the string
# table being recompacted around it.
$ws.Range("D10").Value = "Clear unused reference in project Draw area"

# Insert a new "Status" column before the existing "Description" column
# (old column D -> becomes column E, data/format shift right).
$ws.Columns("D").Insert()

# Give the new column roughly the same width as column C (closest the
# host's pixel-grid column-width quantization allows).
$ws.Columns("D").ColumnWidth = 16

# Header for the new column.
$ws.Range("D2").Value = "Status"

# Fill in the status values. Order matches how they were actually typed in
# (row 5 was revisited/updated last), which also matches the order the
# distinct values were introduced into the workbook.
$ws.Range("D4").Value = "Done"
$ws.Range("D6").Value = "Done"
$ws.Range("D7").Value = "Done"
$ws.Range("D8").Value = "In progress"
$ws.Range("D9").Value = "In progress"
$ws.Range("D10").Value = "In progress"
$ws.Range("D5").Value = "In Analysis"
